# Reorder the "Recorded By" (column G) comma-separated list so that
# "System" (exact case) always comes first, preserving the relative
# order of the remaining names - e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, system, System" -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            }
        }

        if ($hasSystem) {
            $newParts = New-Object System.Collections.ArrayList
            [void]$newParts.Add("System")
            foreach ($p in $parts) {
                if (-not $p.Equals("System")) {
                    [void]$newParts.Add($p)
                }
            }
            $newVal = $newParts -join ", "

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
